$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "antecedents" (col A) and "consequents" (col B) columns encode
# "feature:value" pairs separated by " ∪ ". Rename the separator from
# ":" to "=" for every data row (header row 1 is left untouched).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value()
        if ($val -ne $null) {
            $cell.Value = $val.Replace(":", "=")
        }
    }
}
